$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the Korean translation text for "Maximum size of monument (minimum 3)"
# to the more readable "기념비의 비석 최대 개수 (최소 3)"
$ws.Range("F8").Value = "기념비의 비석 최대 개수 (최소 3)"
